# Update countries & provincias Spain
# Applies the data refresh captured in the diff:
#  - Updated "last updated" timestamp
#  - Peru, Pakistan, Tailandia, Mongolia, Butan, Curazao, Santa Lucia stats updated
#  - Islas Turcas y Caicos overtook Tanzania in total cases, so the two rows swap places

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "data actualizada" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 06:39"

# Row 6 - Peru
$ws.Range("B6").Value = 3769523
$ws.Range("C6").Value = 3415
$ws.Range("D6").Value = 2901908
$ws.Range("E6").Value = 801155

# Row 19 - Pakistan
$ws.Range("B19").Value = 296590
$ws.Range("C19").Value = 441
$ws.Range("D19").Value = 281459
$ws.Range("E19").Value = 8813
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 6318

# Row 124 - Tailandia
$ws.Range("B124").Value = 3425
$ws.Range("C124").Value = 8
$ws.Range("E124").Value = 93

# Rows 172/173 - Islas Turcas y Caicos overtakes Tanzania (rows swap country + stats)
$ws.Range("A172").Value = "Islas Turcas y Caicos"
$ws.Range("B172").Value = 538
$ws.Range("C172").Value = 30
$ws.Range("D172").Value = 218
$ws.Range("E172").Value = 316
$ws.Range("H172").Value = 4

$ws.Range("A173").Value = "Tanzania"
$ws.Range("B173").Value = 509
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 183
$ws.Range("E173").Value = 305
$ws.Range("H173").Value = 21

# Row 183 - Mongolia
$ws.Range("B183").Value = 306
$ws.Range("C183").Value = 2
$ws.Range("D183").Value = 296
$ws.Range("E183").Value = 10

# Row 186 - Butan
$ws.Range("B186").Value = 227
$ws.Range("C186").Value = 2
$ws.Range("E186").Value = 85

# Row 196 - Curazao
$ws.Range("B196").Value = 71
$ws.Range("C196").Value = 2
$ws.Range("E196").Value = 35

# Row 204 - Santa Lucia
$ws.Range("D204").Value = 26
$ws.Range("E204").Value = 0
